{"js": "// 1. Update the report date from 2022-08-21 to 2023-08-21.\nconst dateResults = context.document.body.search(\"2022-08-21\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nfor (const r of dateResults.items) {\n  r.insertText(\"2023-08-21\", Word.InsertLocation.replace);\n}\n\n// 2. Left/start-align both tables in the document (adds <w:jc w:val=\"start\"/> to tblPr).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nfor (const t of tables.items) {\n  t.alignment = \"start\";\n}\n\n// 3. Re-tint a handful of the pandoc syntax-highlighting character styles.\nconst styles = context.document.getStyles();\n\nconst constantTok = styles.getByNameOrNullObject(\"ConstantTok\");\nconstantTok.load(\"font\");\nconst specialCharTok = styles.getByNameOrNullObject(\"SpecialCharTok\");\nspecialCharTok.load(\"font\");\nconst functionTok = styles.getByNameOrNullObject(\"FunctionTok\");\nfunctionTok.load(\"font\");\nconst attributeTok = styles.getByNameOrNullObject(\"AttributeTok\");\nattributeTok.load(\"font\");\nawait context.sync();\n\n// ConstantTok: 000000 -> 8f5902\nconstantTok.font.color = \"#8f5902\";\n\n// SpecialCharTok: 000000 -> ce5c00, add bold\nspecialCharTok.font.color = \"#ce5c00\";\nspecialCharTok.font.bold = true;\n\n// FunctionTok: 000000 -> 204a87, add bold\nfunctionTok.font.color = \"#204a87\";\nfunctionTok.font.bold = true;\n\n// AttributeTok: c4a000 -> 204a87\nattributeTok.font.color = \"#204a87\";\n\nawait context.sync();\n", "ps1": "# Converts a \"RRGGBB\" hex string into the BGR-packed long that the Word\n# object model's Font.Color / WdColor properties expect.\nfunction ConvertTo-WordColor([string]$hex) {\n    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)\n    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)\n    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)\n    return ($b * 65536) + ($g * 256) + $r\n}\n\n$d = $word.ActiveDocument\n\n# 1. Update the report date from 2022-08-21 to 2023-08-21.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"Date\") {\n        $p.Range.Text = \"2023-08-21\"\n    }\n}\n\n# 2. Left/start-align both tables in the document (adds <w:jc w:val=\"start\"/> to tblPr).\nforeach ($t in $d.Tables) {\n    $t.Alignment = \"start\"\n}\n\n# 3. Re-tint a handful of the pandoc syntax-highlighting character styles.\n\n# ConstantTok: 000000 -> 8f5902\n$constantTok = $d.Styles(\"ConstantTok\")\n$constantTok.Font.Color = ConvertTo-WordColor \"8f5902\"\n\n# SpecialCharTok: 000000 -> ce5c00, add bold\n$specialCharTok = $d.Styles(\"SpecialCharTok\")\n$specialCharTok.Font.Color = ConvertTo-WordColor \"ce5c00\"\n$specialCharTok.Font.Bold = $true\n\n# FunctionTok: 000000 -> 204a87, add bold\n$functionTok = $d.Styles(\"FunctionTok\")\n$functionTok.Font.Color = ConvertTo-WordColor \"204a87\"\n$functionTok.Font.Bold = $true\n\n# AttributeTok: c4a000 -> 204a87\n$attributeTok = $d.Styles(\"AttributeTok\")\n$attributeTok.Font.Color = ConvertTo-WordColor \"204a87\"\n"}
